{"js": "// Update the Cycle/date references in the \"Overview\" paragraph:\n//   \"Cycle 1 programs \"            -> \"Cycle 2 and 3 programs \"\n//   \" and Cycle 2 proposals \"      -> \" and Cycle 4 proposals \"\n//   \"in January\"                   -> \"in October\"\n\nasync function replaceFirst(searchText, replacement) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceFirst(\"Cycle 1 programs \", \"Cycle 2 and 3 programs \");\nawait replaceFirst(\" and Cycle 2 proposals \", \" and Cycle 4 proposals \");\nawait replaceFirst(\"in January\", \"in October\");\n", "ps1": "# Update the Cycle/date references in the \"Overview\" paragraph:\n#   \"Cycle 1 programs \"            -> \"Cycle 2 and 3 programs \"\n#   \" and Cycle 2 proposals \"      -> \" and Cycle 4 proposals \"\n#   \"in January\"                   -> \"in October\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        \"wdFindStop\", # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        \"wdReplaceOne\" # Replace\n    )\n}\n\nReplace-FirstMatch \"Cycle 1 programs \" \"Cycle 2 and 3 programs \"\nReplace-FirstMatch \" and Cycle 2 proposals \" \" and Cycle 4 proposals \"\nReplace-FirstMatch \"in January\" \"in October\"\n"}
